$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append a trailing period to the six "answer" cells across rows 3-6 ---
# (D3/F3 = Q3 Parrot answers, D4/F4 = Q4 Shy answers, D5/F5 = Q5 MovieStar
#  answers, D6/F6 = Q6 50years answers)
foreach ($addr in @("D3", "F3", "D4", "F4", "D5", "F5", "D6", "F6")) {
    $cell = $ws.Range($addr)
    $cell.Value = $cell.Value2 + "."
}

# --- Column widths: widen the long-answer columns D and F ---
$ws.Columns.Item(4).ColumnWidth = 90.5
$ws.Columns.Item(6).ColumnWidth = 97.83

# --- View: scroll right so column D is the left-most visible column, and
#     move the active selection to F9 ---
$ws.Range("F9").Select()
$excel.ActiveWindow.ScrollColumn = 4
